$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($range, [string]$text)
    $range.NumberFormat = "@"
    $range.Value = $text
    $range.Style = "Normal"
}

# Rows 2-25: update Price (D) and Volume(1h) (E) values in place
$ws.Range("D2").Value = "27.703.57"
$ws.Range("E2").Value = "  -0.39%  "
$ws.Range("D3").Value = "1.866.02"
$ws.Range("E3").Value = "  -0.79%  "
Set-TextValue $ws.Range("D4") "1.011"
$ws.Range("E4").Value = "  +0.88%  "
Set-TextValue $ws.Range("D5") "333.07"
$ws.Range("E5").Value = "  +0.22%  "
Set-TextValue $ws.Range("D6") "1.010"
$ws.Range("E6").Value = "  +0.67%  "
Set-TextValue $ws.Range("D7") "0.4660"
$ws.Range("E7").Value = "  -1.40%  "
Set-TextValue $ws.Range("D8") "0.3892"
$ws.Range("E8").Value = "  -1.63%  "
Set-TextValue $ws.Range("D9") "46.32"
$ws.Range("E9").Value = "  -3.31%  "
Set-TextValue $ws.Range("D10") "0.07967"
$ws.Range("E10").Value = "  -1.36%  "
Set-TextValue $ws.Range("D11") "0.9986"
$ws.Range("E11").Value = "  -3.21%  "
Set-TextValue $ws.Range("D12") "21.52"
$ws.Range("E12").Value = "  -2.99%  "
$ws.Range("D13").Value = "1.866.84"
$ws.Range("E13").Value = "  -0.97%  "
Set-TextValue $ws.Range("D14") "5.987"
$ws.Range("E14").Value = "  +0.10%  "
Set-TextValue $ws.Range("D15") "7.173"
$ws.Range("E15").Value = "  +0.34%  "
Set-TextValue $ws.Range("D16") "1.012"
$ws.Range("E16").Value = "  +0.72%  "
Set-TextValue $ws.Range("D17") "88.02"
$ws.Range("E17").Value = "  +0.85%  "
Set-TextValue $ws.Range("D18") "0.06720"
$ws.Range("E18").Value = "  +0.62%  "
Set-TextValue $ws.Range("D19") "0.00001045"
$ws.Range("E19").Value = "  -0.63%  "
Set-TextValue $ws.Range("D20") "16.92"
$ws.Range("E20").Value = "  -2.30%  "
Set-TextValue $ws.Range("D21") "1.011"
$ws.Range("E21").Value = "  +0.80%  "
$ws.Range("D22").Value = "27.683.52"
$ws.Range("E22").Value = "  -0.47%  "
Set-TextValue $ws.Range("D23") "5.462"
$ws.Range("E23").Value = "  -1.39%  "
Set-TextValue $ws.Range("D24") "10.88"
$ws.Range("E24").Value = "  -1.29%  "
Set-TextValue $ws.Range("D25") "2.323"
$ws.Range("E25").Value = "  +0.89%  "

# Rows 26-51: coin list shifted up by one (WrappedliquidstakedEther2.0 dropped out of list,
# EOS newly appears at the bottom); update Coin (B), Link (C), Price (D), Volume(1h) (E)
$ws.Range("B26").Value = "Monero"
$ws.Range("C26").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
Set-TextValue $ws.Range("D26") "157.84"
$ws.Range("E26").Value = "  -1.01%  "
$ws.Range("B27").Value = "EthereumClassic"
$ws.Range("C27").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
Set-TextValue $ws.Range("D27") "19.70"
$ws.Range("E27").Value = "  -2.69%  "
$ws.Range("B28").Value = "LidoDAOToken"
$ws.Range("C28").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
Set-TextValue $ws.Range("D28") "2.113"
$ws.Range("E28").Value = "  -0.03%  "
$ws.Range("B29").Value = "InternetComputer(DFINITY)"
$ws.Range("C29").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
Set-TextValue $ws.Range("D29") "5.370"
$ws.Range("E29").Value = "  -4.47%  "
$ws.Range("B30").Value = "BitcoinCash"
$ws.Range("C30").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
Set-TextValue $ws.Range("D30") "121.39"
$ws.Range("E30").Value = "  -0.93%  "
$ws.Range("B31").Value = "ImmutableX"
$ws.Range("C31").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
Set-TextValue $ws.Range("D31") "0.9714"
$ws.Range("E31").Value = "  -1.73%  "
$ws.Range("B32").Value = "Stellar"
$ws.Range("C32").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
Set-TextValue $ws.Range("D32") "0.09445"
$ws.Range("E32").Value = "  -1.17%  "
$ws.Range("B33").Value = "HuobiToken"
$ws.Range("C33").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
Set-TextValue $ws.Range("D33") "3.638"
$ws.Range("E33").Value = "  +1.31%  "
$ws.Range("B34").Value = "Filecoin"
$ws.Range("C34").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
Set-TextValue $ws.Range("D34") "5.299"
$ws.Range("E34").Value = "  -1.45%  "
$ws.Range("B35").Value = "ARBITRUM"
$ws.Range("C35").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
Set-TextValue $ws.Range("D35") "1.329"
$ws.Range("E35").Value = "  -8.47%  "
$ws.Range("B36").Value = "Hedera"
$ws.Range("C36").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
Set-TextValue $ws.Range("D36") "0.06017"
$ws.Range("E36").Value = "  -1.94%  "
$ws.Range("B37").Value = "VeChain"
$ws.Range("C37").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
Set-TextValue $ws.Range("D37") "0.02211"
$ws.Range("E37").Value = "  -2.24%  "
$ws.Range("B38").Value = "TrustWalletToken"
$ws.Range("C38").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
Set-TextValue $ws.Range("D38") "1.197"
$ws.Range("E38").Value = "  -3.28%  "
$ws.Range("B39").Value = "FraxShare"
$ws.Range("C39").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
Set-TextValue $ws.Range("D39") "8.141"
$ws.Range("E39").Value = "  -0.51%  "
$ws.Range("B40").Value = "Frax"
$ws.Range("C40").Value = "https://coinranking.com/coin/KfWtaeV1W+frax-frax"
Set-TextValue $ws.Range("D40") "1.009"
$ws.Range("E40").Value = "  +0.75%  "
$ws.Range("B41").Value = "TheSandbox"
$ws.Range("C41").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
Set-TextValue $ws.Range("D41") "0.5911"
$ws.Range("E41").Value = "  -2.25%  "
$ws.Range("B42").Value = "Algorand"
$ws.Range("C42").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
Set-TextValue $ws.Range("D42") "0.1880"
$ws.Range("E42").Value = "  -1.49%  "
$ws.Range("B43").Value = "Aptos"
$ws.Range("C43").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
Set-TextValue $ws.Range("D43") "10.21"
$ws.Range("E43").Value = "  -0.93%  "
$ws.Range("B44").Value = "WEMIXTOKEN"
$ws.Range("C44").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
Set-TextValue $ws.Range("D44") "1.253"
$ws.Range("E44").Value = "  -1.43%  "
$ws.Range("B45").Value = "Decentraland"
$ws.Range("C45").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
Set-TextValue $ws.Range("D45") "0.5617"
$ws.Range("E45").Value = "  -2.24%  "
$ws.Range("B46").Value = "EnergySwap"
$ws.Range("C46").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
Set-TextValue $ws.Range("D46") "12.03"
$ws.Range("E46").Value = "  -2.08%  "
$ws.Range("B47").Value = "NEARProtocol"
$ws.Range("C47").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
Set-TextValue $ws.Range("D47") "1.915"
$ws.Range("E47").Value = "  -1.89%  "
$ws.Range("B48").Value = "PancakeSwap"
$ws.Range("C48").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
Set-TextValue $ws.Range("D48") "3.297"
$ws.Range("E48").Value = "  -2.71%  "
$ws.Range("B49").Value = "Cronos"
$ws.Range("C49").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
Set-TextValue $ws.Range("D49") "0.06775"
$ws.Range("E49").Value = "  -1.82%  "
$ws.Range("B50").Value = "Quant"
$ws.Range("C50").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
Set-TextValue $ws.Range("D50") "112.06"
$ws.Range("E50").Value = "  -1.97%  "
$ws.Range("B51").Value = "EOS"
$ws.Range("C51").Value = "https://coinranking.com/coin/iAzbfXiBBKkR6+eos-eos"
Set-TextValue $ws.Range("D51") "1.062"
$ws.Range("E51").Value = "  -1.11%  "
